# Update "想去人数" (F column) counts and one "已售罄" (sold out) marker
# across the 展览 (Exhibition), 演出 (Performance) and 全部类型 (All types)
# sheets, per the source data refresh.

$wb = $excel.ActiveWorkbook

$wsExpo = $wb.Worksheets.Item("展览")
$wsShow = $wb.Worksheets.Item("演出")
$wsAll  = $wb.Worksheets.Item("全部类型")

# ---- 展览 (sheet1) ----
$wsExpo.Range("F2").Value = 1374
$wsExpo.Range("F3").Value = 1690
$wsExpo.Range("F4").Value = 905
$wsExpo.Range("F7").Value = 681
$wsExpo.Range("F11").Value = 2480
$wsExpo.Range("F13").Value = 1518
$wsExpo.Range("G14").Value = "已售罄"
$wsExpo.Range("F17").Value = 796
$wsExpo.Range("F18").Value = 86
$wsExpo.Range("F19").Value = 315
$wsExpo.Range("F24").Value = 5125
$wsExpo.Range("F26").Value = 580
$wsExpo.Range("F28").Value = 162
$wsExpo.Range("F30").Value = 228
$wsExpo.Range("F31").Value = 226
$wsExpo.Range("F34").Value = 750
$wsExpo.Range("F37").Value = 239
$wsExpo.Range("F38").Value = 394
$wsExpo.Range("F39").Value = 1088
$wsExpo.Range("F42").Value = 177
$wsExpo.Range("F44").Value = 57

# ---- 演出 (sheet2) ----
$wsShow.Range("F6").Value = 8
$wsShow.Range("F8").Value = 1

# ---- 全部类型 (sheet4) ----
$wsAll.Range("F2").Value = 1374
$wsAll.Range("F5").Value = 1690
$wsAll.Range("F6").Value = 905
$wsAll.Range("F11").Value = 681
$wsAll.Range("F13").Value = 8
$wsAll.Range("F15").Value = 1
$wsAll.Range("F17").Value = 2480
$wsAll.Range("F19").Value = 1518
$wsAll.Range("G20").Value = "已售罄"
$wsAll.Range("F24").Value = 796
$wsAll.Range("F25").Value = 86
$wsAll.Range("F26").Value = 315
$wsAll.Range("F29").Value = 5126
$wsAll.Range("F31").Value = 580
$wsAll.Range("F33").Value = 162
$wsAll.Range("F35").Value = 228
$wsAll.Range("F36").Value = 226
$wsAll.Range("F39").Value = 750
$wsAll.Range("F41").Value = 394
$wsAll.Range("F42").Value = 1088
$wsAll.Range("F44").Value = 177
$wsAll.Range("F46").Value = 57
